# Update the "Manual/snapshot_graphs.xlsx" tracker for the 27-Feb update:
# the newest weekly data point (row 8) moves from 15-Feb-22 to 23-Feb-22
# and its mobility percentages are revised.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 8: date + the five mobility indicator percentages.
$ws.Range("A8").Value2 = 44615      # 23-Feb-22 (serial date)
$ws.Range("B8").Value2 = 0.04       # Retail Mobility
$ws.Range("C8").Value2 = 0.17       # Transit Mobility
$ws.Range("D8").Value2 = 0.12       # Workplace Mobility
$ws.Range("E8").Value2 = 0.48       # Grocery Mobility
$ws.Range("F8").Value2 = 0.08       # Residential Mobility

# Move the active selection to C5 (also clears the old D8 selection /
# scrolled-to-row-8 view state left over from editing that row).
$ws.Range("C5").Select()

$wb.Save()
